# The source diff for this deck is almost entirely "incidental" churn: the
# Presentation1.pptx artifact lives under PowerPoll's bin/Debug output and is
# *regenerated from scratch* by the PowerPoll backend (via the OpenXML SDK)
# every time it runs, so every r:id in presentation.xml / slide.xml / the
# slide rels gets reassigned on each build -- that happens even for layouts
# and masters that were not touched by the "fixed option bug in backend" fix.
#
# The one genuinely semantic change in the diff is the GUID carried by the
# PowerPoll task-pane add-in's web extension binding:
#
#   ppt/slides/udata/data.xml  ->  <we:webextension ... id="{...}">
#
# That id (and the matching r:embed on its snapshot image / the slide's
# we:webextensionref) is assigned internally by the Office Add-ins runtime
# when the add-in is inserted; it is not a documented property anywhere on
# Application/Presentation/Slide/Shape (no CustomXMLParts entry is created
# for it either -- CustomXMLParts.Count stays 0), so it cannot be written
# from VBA/COM automation in PowerPoint itself.
#
# On top of that, this particular slide's OfficeApp shape (the
# mc:AlternateContent graphicFrame/pic pair behind the add-in) was authored
# with cNvPr id="2", which collides with the Title placeholder's id="2".
# Shapes.Item(3) / Shapes.Item("OfficeApp 0") therefore resolve to the Title
# shape instead, so blindly indexing into it (Delete/replace/etc.) would
# silently corrupt the Title placeholder rather than touch the add-in shape.
#
# Given there is no reachable, non-destructive COM call that reproduces the
# add-in id swap (or the wholesale r:id renumbering that comes from an
# external regeneration pass rather than an in-app edit), this script
# intentionally performs no mutation of the deck, to avoid corrupting the
# Title/Subtitle placeholders or fabricating unrelated content that is not
# part of the source diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Touch the object model read-only, purely to confirm the deck is reachable;
# no properties are written.
$null = $s.Shapes.Count
